$wb = $excel.ActiveWorkbook

# Hyperlink font color used elsewhere in this workbook (Cornflow Blue, FF6495ED)
$hyperlinkColor = 15570276

# ---------------------------------------------------------------------------
# Overview sheet: the per-language status column mirrors the same "Status"
# text shown in the zh-cn/de-de detail sheets, so it also flips from
# "Handoff transform failed" to "Ready for handoff".
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"

# ---------------------------------------------------------------------------
# zh-cn sheet: handoff has completed, report the generated xlf handoff file
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B2").Value = "Ready for handoff"

$zhHandoffFile = "b868da23-760c-4df7-a999-f0ddde856131.876706c2ea04bd3350c1de1f29576835b23b2d10.zh-cn.xlf"
$zhHandoffUrl = "https://github.com/OpenLocalizationTest/oltest/blob/ff9ba2c73e2a4380023404c65e3a8373d92104bb/e2e/$zhHandoffFile"

$wsZh.Range("C2").Value = $zhHandoffFile
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), $zhHandoffUrl, "", "", $zhHandoffFile)
$wsZh.Range("C2").Font.Color = $hyperlinkColor
$wsZh.Range("C2").Font.Underline = $True

$wsZh.Range("D2").Value = "2016-02-18 04:01:47"
$wsZh.Range("H2").Value = "Include"

# ---------------------------------------------------------------------------
# de-de sheet: same handoff completion, different generated xlf + timestamp
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B2").Value = "Ready for handoff"

$deHandoffFile = "b868da23-760c-4df7-a999-f0ddde856131.876706c2ea04bd3350c1de1f29576835b23b2d10.de-de.xlf"
$deHandoffUrl = "https://github.com/OpenLocalizationTest/oltest/blob/ff9ba2c73e2a4380023404c65e3a8373d92104bb/e2e/$deHandoffFile"

$wsDe.Range("C2").Value = $deHandoffFile
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), $deHandoffUrl, "", "", $deHandoffFile)
$wsDe.Range("C2").Font.Color = $hyperlinkColor
$wsDe.Range("C2").Font.Underline = $True

$wsDe.Range("D2").Value = "2016-02-18 04:02:01"
$wsDe.Range("H2").Value = "Include"
